$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before FR (this shifts the existing FR -> FS and FS -> FT,
# carrying over the "nom" / "url_produit" headers and all their data).
$ws.Range("FR:FR").Insert()

# New column header: latest snapshot timestamp
$ws.Range("FR1").Value = "2026-02-04 23:16:28"

# The last known price (previously in FQ) is carried forward into the new FR
# column for every product row that already had a price tracked through FQ.
$lastRow = 208
for ($r = 2; $r -le $lastRow; $r++) {
    $prevCell = $ws.Cells.Item($r, 173)   # column FQ
    $newCell = $ws.Cells.Item($r, 174)    # column FR
    $prevValue = $prevCell.Value2
    if ($prevValue -ne $null -and $prevValue -ne "") {
        $newCell.Value = $prevValue
    }
}
